# Fill out the MFix product card template with real item data
# (params taken from the Item), per commit message:
# "1. Card is filled out (with params from Item )"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: bold centered title -------------------------------------------------
# Two-line title needs wrapping + a taller row.
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 50
$ws.Range("B4").Value = "Саморезы гипс/металл`t3.5x25"

# --- Row 5: Marking | YZP ---------------------------------------------------
$ws.Range("B5").Value = "Marking"
# C5 becomes a real (bold, centered) value cell, matching the style already
# used for the other centered value cells (e.g. C6) instead of the plain
# untouched placeholder style it had before.
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "YZP"

# --- Row 6: РАЗМЕР/Size | 3.5x25 --------------------------------------------
$ws.Range("B6").Value = "РАЗМЕР/Size"
$ws.Range("C6").Value = "3.5x25"

# --- Row 7: spacer row, label + value cleared back to plain formatting -----
$ws.Range("B10").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").ClearContents()

$ws.Range("D10").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").ClearContents()

# --- Row 8: Кол-во в упак/шт. | 1000 | Шт / PCS -----------------------------
$ws.Range("B8").Value = "Кол-во в упак/шт."
# "1000" must stay a text value (matches the rest of the card, which stores
# every value as shared-string text) instead of being auto-detected as a
# number, so force Text format for the entry then restore the original
# (General / same style as the sibling value cell C6) formatting.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "1000"
$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D8").Value = "Шт / PCS"

# --- Row 9: Вес упак Кг/Kgs | (blank) | Кг/Kgs ------------------------------
$ws.Range("B9").Value = "Вес упак Кг/Kgs"
$ws.Range("D9").Value = "Кг/Kgs"

# --- Row 10: (blank) | Сделано в КНР | (blank) ------------------------------
$ws.Range("C10").Value = "Сделано в КНР"

# --- Row 11: ORDER: | 2155695PL ---------------------------------------------
$ws.Range("B11").Value = "ORDER:"
$ws.Range("C11").Value = "2155695PL"
